# Update values in Sheet1 (result_data_KNN) per commit "Update Name of Algo"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A3").Value = -21.76
$ws.Range("D3").Value = -7.768000000000001
$ws.Range("A21").Value = -20.469
$ws.Range("A23").Value = -20.428
$ws.Range("D24").Value = -7.734999999999999
$ws.Range("A25").Value = -21.664
$ws.Range("B27").Value = 6.152000000000001
$ws.Range("B31").Value = 6.216
$ws.Range("B39").Value = 7.932
$ws.Range("B48").Value = 5.323
$ws.Range("B51").Value = 6.226
$ws.Range("B52").Value = 5.813000000000001
$ws.Range("A53").Value = -21.945
$ws.Range("B55").Value = 4.671
$ws.Range("B56").Value = 5.003
$ws.Range("A57").Value = -21.303
$ws.Range("B57").Value = 6.432
$ws.Range("D57").Value = -8.231999999999999
$ws.Range("A59").Value = -22.397
$ws.Range("D61").Value = -7.704000000000001
$ws.Range("A69").Value = -21.649
$ws.Range("D70").Value = -7.132
$ws.Range("B73").Value = 6.751
$ws.Range("A79").Value = -21.137
$ws.Range("A83").Value = -22.065
$ws.Range("D86").Value = -8.246
$ws.Range("B89").Value = 5.946
$ws.Range("B90").Value = 5.671
$ws.Range("A93").Value = -21.398
$ws.Range("D98").Value = -8.099
$ws.Range("D100").Value = -8.022000000000002
$ws.Range("D102").Value = -8.285
